# Applies the commit's edit:
# - Rewrites the single shared-string cell's text from a compact one-line
#   Python-dict-literal dump into a pretty-printed, JSON-quoted, multi-line
#   rendering of the same "questions" data.
# - Moves that text from A2 up into A1 (the row-2 cell is removed).
# - Strips the bold+bordered/centered style that used to live on A1, so the
#   cell falls back to the workbook's default (unstyled) formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = @'
questions = [
    {
        "title": "You are a network engineer at a multinational corporation. You're designing a communication system for a department needing real-time video data streaming. The system must transmit video data smoothly and without delay.Which protocol should you implement for this system?",
        "ques_type": 2,
        "options": [
            "UDP (User Datagram Protocol)",
            "ICMP (Internet Control Message Protocol)",
            "SNMP (Simple Network Management Protocol)",
            "FTP (File Transfer Protocol)"
        ],
        "score": "UDP (User Datagram Protocol)"
    },
    {
        "title": "You are a network engineer at an expanding multinational corporation. You need to redesign the existing LAN (local area network) to accommodate more employees and increased network resource demand. The objective is a scalable, efficient LAN topology design that supports high-speed data transfer, ensures network reliability, and provides seamless connectivity across various departments and locations.What should be your primary consideration in designing the LAN topology?",
        "ques_type": 2,
        "options": [
            "Establishing a hierarchical design with core, distribution, and access layers.",
            "Implementing redundant links and devices. ",
            "Considering bandwidth requirements and choosing suitable network equipment.",
            "Creating virtual local area networks. "
        ],
        "score": "Considering bandwidth requirements and choosing suitable network equipment."
    },
    {
        "title": "You are a network engineer at a large financial institution responsible for network security. You've recently implemented a deep packet inspection (DPI) firewall, which has just alerted you to unauthorized file transfers of sensitive company data to an external email address.What should be your immediate action?",
        "ques_type": 2,
        "options": [
            "Cut off access to the external email address and terminate the employee's network access.",
            "Implement extra security measures such as two-factor authentication for accessing sensitive data.",
            "Adjust the firewall to automatically block any file transfers containing sensitive data.",
            "Capture and scrutinize the packet data associated with the file transfer for investigation."
        ],
        "score": "Capture and scrutinize the packet data associated with the file transfer for investigation."
    },
    {
        "title": "You are a network engineer at a large telecommunications company. You are responsible for the performance of the company's network infrastructure. Recently, multiple users have reported slow internet speeds and intermittent connectivity issues. You suspect congestion on one of the network switches.Which action should you take next?",
        "ques_type": 2,
        "options": [
            "Inspect the physical connections of the network switch.",
            "Reset the network switch.",
            "Upgrade the network switch to a model with higher capacity.",
            "Review the Simple Network Management Protocol (SNMP) data from the switch."
        ],
        "score": "Review the Simple Network Management Protocol (SNMP) data from the switch."
    }
]
'@

# A1 used to hold a plain 0 with a bold/bordered/centered style ("s=1");
# drop that formatting so the cell reverts to the default style.
$ws.Range("A1").ClearFormats()

# A2 held the shared-string text; put the (reformatted) text in A1 instead.
$ws.Range("A1").Value = $newText

# Row 2 (which used to carry the text) is no longer needed.
$ws.Rows(2).Delete()

# The multi-line text auto-expanded the row height; auto-fit it back down
# so no explicit row height is persisted.
$ws.Rows(1).AutoFit()

Write-Host "done"
